$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: swap the "id_carga" figures and the explanatory text
$ws.Range("B3").Value = "id_carga: 20, id_buque:1"
$ws.Range("C3").Value = "No se deben actualizar las tablas alojamiento_bodega ni Cargar_maritima, ya que dicho buque no posee la carga descrita."

# New row 4
$ws.Range("A4").Value = "Requerimiento 11"
$ws.Range("B4").Value = "id_carga: 2, id_buque:2"
$ws.Range("C4").Value = "Se actualiza alguna de las bodegas. En la tabla alojamiento_bodega. Se elimina la fila correspondiente de la tabla Cargar_maritima."

# New row 5
$ws.Range("A5").Value = "Requerimiento 11"
$ws.Range("B5").Value = "id_carga: 3, id_buque:3"
$ws.Range("C5").Value = "Se actualiza alguna de las bodegas. En la tabla alojamiento_bodega. Se elimina la fila correspondiente de la tabla Cargar_maritima."

# Match formatting used on the existing rows (wrap text + taller rows)
$ws.Range("C4:C5").WrapText = $true
$ws.Range("A4:C5").RowHeight = 60

# Update the saved selection to D3
$ws.Range("D3").Select()
